# Fruta / hortaliza, semanal
# Re-sorts/updates the weekly price rows (2-12) for Membrillo @ Vega Monumental
# Concepcion so each row carries the Fecha/Calidad/Volumen/Precio*/Unidad/Origen
# combination for its (now reordered) week. Columns A,B,C,E,F,G,H,I,J,K are
# identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44363; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 15 kilos empedrada";   R="Región de O'Higgins"; S=633; T=15 },
    @{ Row=3;  D=44425; L="Primera"; M=100; N=12000; O=13000; P=12500; Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=694; T=18 },
    @{ Row=4;  D=44299; L="Primera"; M=100; N=10000; O=11000; P=10500; Q="`$/caja 18 kilos granel";      R="Región del Maule";    S=583; T=18 },
    @{ Row=5;  D=44299; L="Segunda"; M=50;  N=9000;  O=9000;  P=9000;  Q="`$/caja 18 kilos granel";      R="Región del Maule";    S=500; T=18 },
    @{ Row=6;  D=44272; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 15 kilos granel";      R="Región de O'Higgins"; S=633; T=15 },
    @{ Row=7;  D=44272; L="Segunda"; M=50;  N=8000;  O=8000;  P=8000;  Q="`$/caja 15 kilos granel";      R="Región de O'Higgins"; S=533; T=15 },
    @{ Row=8;  D=44307; L="Primera"; M=50;  N=10000; O=10000; P=10000; Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=556; T=18 },
    @{ Row=9;  D=44307; L="Segunda"; M=50;  N=8000;  O=8000;  P=8000;  Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=444; T=18 },
    @{ Row=10; D=44698; L="Primera"; M=50;  N=10000; O=10000; P=10000; Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=556; T=18 },
    @{ Row=11; D=44316; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=528; T=18 },
    @{ Row=12; D=44358; L="Primera"; M=100; N=11000; O=12000; P=11500; Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=639; T=18 }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("N$r").Value = $rec.N
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R
    $ws.Range("S$r").Value = $rec.S
    $ws.Range("T$r").Value = $rec.T
}
